$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jogadores")

# --- Reset/reapply the AutoFilter while the data still ends at row 293,
#     so the generated ref stays "A1:L293" and is not auto-expanded later
#     once the new rows are appended below. ---
$ws.AutoFilterMode = $false
$ws.Range("A1:L293").AutoFilter(1)

# --- Update F261 value from 1 to 2 ---
$ws.Cells.Item(261, 6).Value = 2

# --- Add new player rows 294-315 (new week of games) ---
$ws.Cells.Item(294, 1).Value = "Cabeleira"
$ws.Cells.Item(294, 3).Value = 5
$ws.Cells.Item(294, 4).Value = 2
$ws.Cells.Item(294, 5).Value = 2
$ws.Cells.Item(294, 6).Value = 2
$ws.Cells.Item(294, 7).Value = 1
$ws.Cells.Item(294, 8).Value = 1
$ws.Cells.Item(294, 9).Value = 0
$ws.Cells.Item(294, 10).Value = 0
$ws.Cells.Item(294, 11).Value = 0
$ws.Cells.Item(294, 12).Value = 0

$ws.Cells.Item(295, 1).Value = "Guinha"
$ws.Cells.Item(295, 3).Value = 5
$ws.Cells.Item(295, 4).Value = 2
$ws.Cells.Item(295, 5).Value = 2
$ws.Cells.Item(295, 6).Value = 0
$ws.Cells.Item(295, 7).Value = 1
$ws.Cells.Item(295, 8).Value = 1
$ws.Cells.Item(295, 9).Value = 0
$ws.Cells.Item(295, 10).Value = 0
$ws.Cells.Item(295, 11).Value = 0
$ws.Cells.Item(295, 12).Value = 0

$ws.Cells.Item(296, 1).Value = "Fernando"
$ws.Cells.Item(296, 3).Value = 5
$ws.Cells.Item(296, 4).Value = 2
$ws.Cells.Item(296, 5).Value = 2
$ws.Cells.Item(296, 6).Value = 2
$ws.Cells.Item(296, 7).Value = 1
$ws.Cells.Item(296, 8).Value = 1
$ws.Cells.Item(296, 9).Value = 0
$ws.Cells.Item(296, 10).Value = 0
$ws.Cells.Item(296, 11).Value = 0
$ws.Cells.Item(296, 12).Value = 0

$ws.Cells.Item(297, 1).Value = "Alan"
$ws.Cells.Item(297, 3).Value = 5
$ws.Cells.Item(297, 4).Value = 2
$ws.Cells.Item(297, 5).Value = 2
$ws.Cells.Item(297, 6).Value = 2
$ws.Cells.Item(297, 7).Value = 1
$ws.Cells.Item(297, 8).Value = 1
$ws.Cells.Item(297, 9).Value = 0
$ws.Cells.Item(297, 10).Value = 0
$ws.Cells.Item(297, 11).Value = 0
$ws.Cells.Item(297, 12).Value = 0

$ws.Cells.Item(298, 1).Value = "Leandrinho"
$ws.Cells.Item(298, 3).Value = 5
$ws.Cells.Item(298, 4).Value = 2
$ws.Cells.Item(298, 5).Value = 2
$ws.Cells.Item(298, 6).Value = 4
$ws.Cells.Item(298, 7).Value = 1
$ws.Cells.Item(298, 8).Value = 1
$ws.Cells.Item(298, 9).Value = 0
$ws.Cells.Item(298, 10).Value = 1
$ws.Cells.Item(298, 11).Value = 0
$ws.Cells.Item(298, 12).Value = 0

$ws.Cells.Item(299, 1).Value = "Boneco"
$ws.Cells.Item(299, 3).Value = 2
$ws.Cells.Item(299, 4).Value = 3
$ws.Cells.Item(299, 5).Value = 2
$ws.Cells.Item(299, 6).Value = 2
$ws.Cells.Item(299, 7).Value = 1
$ws.Cells.Item(299, 8).Value = 0
$ws.Cells.Item(299, 9).Value = 0
$ws.Cells.Item(299, 10).Value = 0
$ws.Cells.Item(299, 11).Value = 0
$ws.Cells.Item(299, 12).Value = 0

$ws.Cells.Item(300, 1).Value = "Juscielio"
$ws.Cells.Item(300, 3).Value = 2
$ws.Cells.Item(300, 4).Value = 3
$ws.Cells.Item(300, 5).Value = 2
$ws.Cells.Item(300, 6).Value = 2
$ws.Cells.Item(300, 7).Value = 1
$ws.Cells.Item(300, 8).Value = 0
$ws.Cells.Item(300, 9).Value = 0
$ws.Cells.Item(300, 10).Value = 0
$ws.Cells.Item(300, 11).Value = 0
$ws.Cells.Item(300, 12).Value = 0

$ws.Cells.Item(301, 1).Value = "Marcos"
$ws.Cells.Item(301, 3).Value = 2
$ws.Cells.Item(301, 4).Value = 3
$ws.Cells.Item(301, 5).Value = 2
$ws.Cells.Item(301, 6).Value = 2
$ws.Cells.Item(301, 7).Value = 1
$ws.Cells.Item(301, 8).Value = 0
$ws.Cells.Item(301, 9).Value = 0
$ws.Cells.Item(301, 10).Value = 0
$ws.Cells.Item(301, 11).Value = 0
$ws.Cells.Item(301, 12).Value = 0

$ws.Cells.Item(302, 1).Value = "Corinthiano"
$ws.Cells.Item(302, 3).Value = 2
$ws.Cells.Item(302, 4).Value = 3
$ws.Cells.Item(302, 5).Value = 2
$ws.Cells.Item(302, 6).Value = 0
$ws.Cells.Item(302, 7).Value = 1
$ws.Cells.Item(302, 8).Value = 0
$ws.Cells.Item(302, 9).Value = 0
$ws.Cells.Item(302, 10).Value = 0
$ws.Cells.Item(302, 11).Value = 0
$ws.Cells.Item(302, 12).Value = 0

$ws.Cells.Item(303, 1).Value = "Michel"
$ws.Cells.Item(303, 3).Value = 2
$ws.Cells.Item(303, 4).Value = 3
$ws.Cells.Item(303, 5).Value = 2
$ws.Cells.Item(303, 6).Value = 0
$ws.Cells.Item(303, 7).Value = 1
$ws.Cells.Item(303, 8).Value = 0
$ws.Cells.Item(303, 9).Value = 0
$ws.Cells.Item(303, 10).Value = 0
$ws.Cells.Item(303, 11).Value = 0
$ws.Cells.Item(303, 12).Value = 0

$ws.Cells.Item(304, 1).Value = "Leandrão"
$ws.Cells.Item(304, 3).Value = 1
$ws.Cells.Item(304, 4).Value = 3
$ws.Cells.Item(304, 5).Value = 2
$ws.Cells.Item(304, 6).Value = 0
$ws.Cells.Item(304, 7).Value = 1
$ws.Cells.Item(304, 8).Value = 0
$ws.Cells.Item(304, 9).Value = 0
$ws.Cells.Item(304, 10).Value = 0
$ws.Cells.Item(304, 11).Value = 0
$ws.Cells.Item(304, 12).Value = 0

$ws.Cells.Item(305, 1).Value = "Marcelão"
$ws.Cells.Item(305, 3).Value = 1
$ws.Cells.Item(305, 4).Value = 3
$ws.Cells.Item(305, 5).Value = 2
$ws.Cells.Item(305, 6).Value = 0
$ws.Cells.Item(305, 7).Value = 1
$ws.Cells.Item(305, 8).Value = 0
$ws.Cells.Item(305, 9).Value = 0
$ws.Cells.Item(305, 10).Value = 0
$ws.Cells.Item(305, 11).Value = 0
$ws.Cells.Item(305, 12).Value = 0

$ws.Cells.Item(306, 1).Value = "Ismael"
$ws.Cells.Item(306, 3).Value = 1
$ws.Cells.Item(306, 4).Value = 3
$ws.Cells.Item(306, 5).Value = 2
$ws.Cells.Item(306, 6).Value = 0
$ws.Cells.Item(306, 7).Value = 1
$ws.Cells.Item(306, 8).Value = 0
$ws.Cells.Item(306, 9).Value = 0
$ws.Cells.Item(306, 10).Value = 0
$ws.Cells.Item(306, 11).Value = 0
$ws.Cells.Item(306, 12).Value = 0

$ws.Cells.Item(307, 1).Value = "Coxinha"
$ws.Cells.Item(307, 3).Value = 1
$ws.Cells.Item(307, 4).Value = 3
$ws.Cells.Item(307, 5).Value = 2
$ws.Cells.Item(307, 6).Value = 2
$ws.Cells.Item(307, 7).Value = 1
$ws.Cells.Item(307, 8).Value = 0
$ws.Cells.Item(307, 9).Value = 0
$ws.Cells.Item(307, 10).Value = 0
$ws.Cells.Item(307, 11).Value = 0
$ws.Cells.Item(307, 12).Value = 0

$ws.Cells.Item(308, 1).Value = "Du"
$ws.Cells.Item(308, 3).Value = 1
$ws.Cells.Item(308, 4).Value = 3
$ws.Cells.Item(308, 5).Value = 2
$ws.Cells.Item(308, 6).Value = 2
$ws.Cells.Item(308, 7).Value = 1
$ws.Cells.Item(308, 8).Value = 0
$ws.Cells.Item(308, 9).Value = 0
$ws.Cells.Item(308, 10).Value = 0
$ws.Cells.Item(308, 11).Value = 0
$ws.Cells.Item(308, 12).Value = 0

$ws.Cells.Item(309, 1).Value = "Said"
$ws.Cells.Item(309, 3).Value = 1
$ws.Cells.Item(309, 4).Value = 2
$ws.Cells.Item(309, 5).Value = 3
$ws.Cells.Item(309, 6).Value = 1
$ws.Cells.Item(309, 7).Value = 1
$ws.Cells.Item(309, 8).Value = 0
$ws.Cells.Item(309, 9).Value = 1
$ws.Cells.Item(309, 10).Value = 0
$ws.Cells.Item(309, 11).Value = 0
$ws.Cells.Item(309, 12).Value = 0

$ws.Cells.Item(310, 1).Value = "Leo"
$ws.Cells.Item(310, 3).Value = 1
$ws.Cells.Item(310, 4).Value = 2
$ws.Cells.Item(310, 5).Value = 3
$ws.Cells.Item(310, 6).Value = 0
$ws.Cells.Item(310, 7).Value = 1
$ws.Cells.Item(310, 8).Value = 0
$ws.Cells.Item(310, 9).Value = 1
$ws.Cells.Item(310, 10).Value = 0
$ws.Cells.Item(310, 11).Value = 0
$ws.Cells.Item(310, 12).Value = 0

$ws.Cells.Item(311, 1).Value = "Euler"
$ws.Cells.Item(311, 3).Value = 1
$ws.Cells.Item(311, 4).Value = 2
$ws.Cells.Item(311, 5).Value = 3
$ws.Cells.Item(311, 6).Value = 0
$ws.Cells.Item(311, 7).Value = 1
$ws.Cells.Item(311, 8).Value = 0
$ws.Cells.Item(311, 9).Value = 1
$ws.Cells.Item(311, 10).Value = 0
$ws.Cells.Item(311, 11).Value = 0
$ws.Cells.Item(311, 12).Value = 0

$ws.Cells.Item(312, 1).Value = "Eder"
$ws.Cells.Item(312, 3).Value = 1
$ws.Cells.Item(312, 4).Value = 2
$ws.Cells.Item(312, 5).Value = 3
$ws.Cells.Item(312, 6).Value = 0
$ws.Cells.Item(312, 7).Value = 1
$ws.Cells.Item(312, 8).Value = 0
$ws.Cells.Item(312, 9).Value = 1
$ws.Cells.Item(312, 10).Value = 0
$ws.Cells.Item(312, 11).Value = 0
$ws.Cells.Item(312, 12).Value = 0

$ws.Cells.Item(313, 1).Value = "Vander"
$ws.Cells.Item(313, 3).Value = 1
$ws.Cells.Item(313, 4).Value = 2
$ws.Cells.Item(313, 5).Value = 3
$ws.Cells.Item(313, 6).Value = 1
$ws.Cells.Item(313, 7).Value = 1
$ws.Cells.Item(313, 8).Value = 0
$ws.Cells.Item(313, 9).Value = 1
$ws.Cells.Item(313, 10).Value = 0
$ws.Cells.Item(313, 11).Value = 0
$ws.Cells.Item(313, 12).Value = 0

$ws.Cells.Item(314, 1).Value = "Matheus"
$ws.Cells.Item(314, 3).Value = 5
$ws.Cells.Item(314, 4).Value = 5
$ws.Cells.Item(314, 5).Value = 4
$ws.Cells.Item(314, 6).Value = 0
$ws.Cells.Item(314, 7).Value = 1
$ws.Cells.Item(314, 8).Value = 1
$ws.Cells.Item(314, 9).Value = 0
$ws.Cells.Item(314, 10).Value = 0
$ws.Cells.Item(314, 11).Value = 11
$ws.Cells.Item(314, 12).Value = 1

$ws.Cells.Item(315, 1).Value = "Chelin"
$ws.Cells.Item(315, 3).Value = 4
$ws.Cells.Item(315, 4).Value = 5
$ws.Cells.Item(315, 5).Value = 5
$ws.Cells.Item(315, 6).Value = 0
$ws.Cells.Item(315, 7).Value = 1
$ws.Cells.Item(315, 8).Value = 0
$ws.Cells.Item(315, 9).Value = 1
$ws.Cells.Item(315, 10).Value = 0
$ws.Cells.Item(315, 11).Value = 11
$ws.Cells.Item(315, 12).Value = 0

# --- Update the _xlnm._FilterDatabase defined name range ---
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Jogadores!`$A`$1:`$L`$293"

# --- Restore the frozen header-row pane and update the active selection ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I315").Select()
